# ------------------------------------------------------------------------
# Applies the "Added ToDataSet/ToDataTable extensions" workbook changes:
#   1. Moves the TEST4 table (and its backing data) from A1:C5 to F9:H13
#      on sheet "TEST4".
#   2. Adds a brand-new worksheet "TEST5" (placed after TEST4, becomes the
#      active sheet/tab) containing a small barcode table in A1:C4.
# ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- Style templates we will reuse (format-only copy/paste) -----------
$xlPasteFormats = -4122

# =========================================================================
# STEP 1: Move TEST4's data block from A1:C5 -> F9:H13
# =========================================================================
$ws4 = $wb.Worksheets.Item("TEST4")

$moveMap = @(
    @("A1", "F9"),  @("B1", "G9"),  @("C1", "H9"),
    @("A2", "F10"),                 @("C2", "H10"),
    @("A3", "F11"), @("B3", "G11"), @("C3", "H11"),
    @("A4", "F12"), @("B4", "G12"), @("C4", "H12"),
    @("A5", "F13"), @("B5", "G13"), @("C5", "H13")
)

foreach ($pair in $moveMap) {
    $ws4.Range($pair[0]).Copy($ws4.Range($pair[1]))
}

# Re-anchor the TEST4 table/autofilter onto its new range.
$loTest4 = $ws4.ListObjects.Item("TEST4")
$loTest4.Resize($ws4.Range("F9:H13"))

# Wipe out the old A1:C5 block now that everything has been copied out.
$ws4.Range("A1:C5").Clear()

# Update the view: the old selection pointed at F6; the new one sits on
# the moved header cell, H9.
$ws4.Range("H9").Select() | Out-Null

# =========================================================================
# STEP 2: Add the new "TEST5" worksheet after TEST4
# =========================================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws5.Name = "TEST5"

# -- Header row (bold "TableStyleLight"-ish header look reused from the
#    other sheets: font-bold style already present in the workbook). -----
$ws5.Range("A1").Value2 = "Barcode"
$ws5.Range("B1").Value2 = "Quantity"
$ws5.Range("C1").Value2 = "UpdatedDate"
$ws5.Range("A1:C1").Font.Bold = $true

# -- Data rows. The shared-strings table records *first-write* order, so
#    to land "Barcode2"/"Barcode1"/"Barcode3" at shared-string indexes
#    44/45/46 respectively we must literally write A3 before A2 before A4
#    (even though the rows end up in 2,3,4 order on the sheet). -----------
$ws5.Range("A3").Value2 = "Barcode2"
$ws5.Range("A2").Value2 = "Barcode1"
$ws5.Range("A4").Value2 = "Barcode3"

$ws5.Range("B2").Value2 = 23
$ws5.Range("B3").Value2 = 12
$ws5.Range("B4").Value2 = 2

# Date column: copy the existing short-date cell format (style already
# used elsewhere in the workbook) onto C2:C4, then stamp in the serial
# date values (2017-08-08).
$ws1 = $wb.Worksheets.Item("TEST1")
$ws1.Range("B2").Copy()
$ws5.Range("C2:C4").PasteSpecial($xlPasteFormats)

$ws5.Range("C2").Value2 = 42955
$ws5.Range("C3").Value2 = 42955
$ws5.Range("C4").Value2 = 42955

# Column C width (roughly matches the "best fit" width used elsewhere).
$ws5.Columns.Item(3).ColumnWidth = 12

# Page setup to match the sibling sheets.
$ws5.PageSetup.PaperSize = 9
$ws5.PageSetup.Orientation = 1

# Final selection/active cell for TEST5's view, and make TEST5 the active
# (selected) sheet/tab, matching tabSelected moving from TEST4 to TEST5.
$ws5.Range("H14").Select() | Out-Null
$ws5.Activate()
